$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $origStyle = $range.Style
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = $origStyle
}

Set-TextValue $ws.Range("D2") "27.010.41"
$ws.Range("E2").Value = "  +2.85%  "

Set-TextValue $ws.Range("D3") "1.653.74"
$ws.Range("E3").Value = "  +3.72%  "

$ws.Range("E4").Value = "  +0.22%  "

Set-TextValue $ws.Range("D5") "214.96"
$ws.Range("E5").Value = "  +1.53%  "

Set-TextValue $ws.Range("D6") "0.508"
$ws.Range("E6").Value = "  +0.67%  "

$ws.Range("E7").Value = "  +0.13%  "

Set-TextValue $ws.Range("D8") "0.249"
$ws.Range("E8").Value = "  +1.72%  "

Set-TextValue $ws.Range("D9") "0.0614"
$ws.Range("E9").Value = "  +1.39%  "

Set-TextValue $ws.Range("D10") "19.63"
$ws.Range("E10").Value = "  +3.54%  "

Set-TextValue $ws.Range("D11") "0.0865"
$ws.Range("E11").Value = "  +1.46%  "

$ws.Range("E12").Value = "  +3.67%  "

Set-TextValue $ws.Range("D13") "1.644.89"
$ws.Range("E13").Value = "  +3.06%  "

Set-TextValue $ws.Range("D14") "4.08"
$ws.Range("E14").Value = "  +1.92%  "

$ws.Range("E15").Value = "  +3.19%  "

Set-TextValue $ws.Range("D16") "64.79"
$ws.Range("E16").Value = "  +1.83%  "

Set-TextValue $ws.Range("D17") "27.014.50"
$ws.Range("E17").Value = "  +2.92%  "

Set-TextValue $ws.Range("D18") "237.46"
$ws.Range("E18").Value = "  +3.29%  "

Set-TextValue $ws.Range("D19") "7.85"
$ws.Range("E19").Value = "  +2.53%  "

$ws.Range("E20").Value = "  +1.04%  "

$ws.Range("E21").Value = "  +0.24%  "

$ws.Range("E22").Value = "  +4.64%  "

$ws.Range("E23").Value = "  +3.97%  "

Set-TextValue $ws.Range("D24") "9.27"
$ws.Range("E24").Value = "  +3.90%  "

Set-TextValue $ws.Range("D25") "146.15"
$ws.Range("E25").Value = "  +0.37%  "

$ws.Range("E26").Value = "  +0.13%  "

Set-TextValue $ws.Range("D27") "7.15"
$ws.Range("E27").Value = "  +2.31%  "

$ws.Range("E28").Value = "  +0.83%  "

$ws.Range("E29").Value = "  +3.16%  "

$ws.Range("E30").Value = "  +0.70%  "

$ws.Range("E31").Value = "  +1.49%  "

Set-TextValue $ws.Range("D32") "1.524.66"
$ws.Range("E32").Value = "  +4.15%  "

Set-TextValue $ws.Range("D33") "3.29"
$ws.Range("E33").Value = "  +2.79%  "

$ws.Range("E34").Value = "  +3.11%  "

$ws.Range("E35").Value = "  +8.25%  "

$ws.Range("E36").Value = "  -0.03%  "

Set-TextValue $ws.Range("D37") "0.575"
$ws.Range("E37").Value = "  +1.43%  "

Set-TextValue $ws.Range("D38") "0.885"
$ws.Range("E38").Value = "  +7.94%  "

$ws.Range("E39").Value = "  +2.55%  "

$ws.Range("E40").Value = "  +2.68%  "

$ws.Range("E41").Value = "  +0.19%  "

Set-TextValue $ws.Range("D44") "1.793.56"
$ws.Range("E44").Value = "  +3.54%  "

$ws.Range("E45").Value = "  +2.35%  "

Set-TextValue $ws.Range("D46") "0.918"
$ws.Range("E46").Value = "  -1.23%  "

Set-TextValue $ws.Range("D47") "90.04"
$ws.Range("E47").Value = "  +2.80%  "

$ws.Range("E48").Value = "  +0.86%  "

Set-TextValue $ws.Range("D49") "1.53"
$ws.Range("E49").Value = "  +3.37%  "

$ws.Range("E50").Value = "  +0.71%  "

Set-TextValue $ws.Range("D51") "0.0976"
$ws.Range("E51").Value = "  +3.12%  "

# Rows 42/43: Aave and MXToken swap positions with updated data
$ws.Range("B42").Value = "Aave"
$ws.Range("C42").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
Set-TextValue $ws.Range("D42") "66.13"
$ws.Range("E42").Value = "  +9.40%  "

$ws.Range("B43").Value = "MXToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
Set-TextValue $ws.Range("D43") "2.25"
$ws.Range("E43").Value = "  +3.61%  "

